$p = $ppt.ActivePresentation

# The table on slide 6 ("SOURCES OF FINANCE") had its table design
# changed (PowerPoint Table Design gallery) from the custom "Table_0"
# style to a different built-in table style.
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{1CC19128-4BEB-40E3-A4BC-FE1E8FDE1194}")
    }
}
